# Mock Project submission by Jan
# TestData.xlsx / LoginPage sheet: replace the sample login rows with a
# single real-looking credential pair, and turn the user name (an email
# address) into a clickable mailto hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the only remaining data row: UserName / Password values.
$ws.Range("A2").Value = "vuppal@gmail.com"
$ws.Range("B2").Value = "test123"

# Drop the old extra sample rows (Testing1/Testing2) - only header + one
# data row remain afterwards.
$ws.Rows("3:4").Delete()

# Excel auto-links an email address typed into a cell; reproduce that by
# adding a mailto hyperlink on the UserName cell.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:vuppal@gmail.com")

# Leave the selection where the author last clicked before saving.
$ws.Range("C11").Select()
